$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value in column D looks numeric (e.g. "1.00", "0.590")
# must be forced to Text format first, otherwise Excel auto-converts the
# assigned string into a number (losing the literal text representation),
# same as typing such a value directly into a General-formatted cell.
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D51').NumberFormat = "@"

$ws.Range('D2').Value = '63.533.68'
$ws.Range('E2').Value = '  -6.97%  '

$ws.Range('D3').Value = '3.513.74'
$ws.Range('E3').Value = '  -3.35%  '

$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.28%  '

$ws.Range('D5').Value = '394.52'
$ws.Range('E5').Value = '  -6.23%  '

$ws.Range('D6').Value = '122.43'
$ws.Range('E6').Value = '  -7.17%  '

$ws.Range('D7').Value = '3.503.96'
$ws.Range('E7').Value = '  -3.39%  '

$ws.Range('D8').Value = '0.590'
$ws.Range('E8').Value = '  -9.79%  '

$ws.Range('E9').Value = '  +0.08%  '

$ws.Range('D10').Value = '0.681'
$ws.Range('E10').Value = '  -12.40%  '

$ws.Range('D11').Value = '0.153'
$ws.Range('E11').Value = '  -15.96%  '

$ws.Range('D12').Value = '0.0000335'
$ws.Range('E12').Value = '  -6.43%  '

$ws.Range('D13').Value = '39.01'
$ws.Range('E13').Value = '  -9.12%  '

$ws.Range('D14').Value = '4.053.34'
$ws.Range('E14').Value = '  -2.70%  '

$ws.Range('D15').Value = '9.22'
$ws.Range('E15').Value = '  -8.13%  '

$ws.Range('E16').Value = '  -3.29%  '

$ws.Range('D17').Value = '3.499.56'
$ws.Range('E17').Value = '  -4.09%  '

$ws.Range('B18').Value = 'Uniswap'
$ws.Range('C18').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D18').Value = '12.75'
$ws.Range('E18').Value = '  +2.14%  '

$ws.Range('B19').Value = 'Chainlink'
$ws.Range('C19').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D19').Value = '18.77'
$ws.Range('E19').Value = '  -8.84%  '

$ws.Range('D20').Value = '63.646.50'
$ws.Range('E20').Value = '  -6.47%  '

$ws.Range('D21').Value = '1.02'
$ws.Range('E21').Value = '  -11.33%  '

$ws.Range('D22').Value = '395.92'
$ws.Range('E22').Value = '  -15.26%  '

$ws.Range('D23').Value = '13.86'
$ws.Range('E23').Value = '  +2.93%  '

$ws.Range('D24').Value = '81.38'
$ws.Range('E24').Value = '  -9.10%  '

$ws.Range('D25').Value = '2.88'
$ws.Range('E25').Value = '  -8.67%  '

$ws.Range('D26').Value = '33.62'
$ws.Range('E26').Value = '  -7.89%  '

$ws.Range('E27').Value = '  +6.95%  '

$ws.Range('D28').Value = '2.99'
$ws.Range('E28').Value = '  -11.64%  '

$ws.Range('D29').Value = '8.84'
$ws.Range('E29').Value = '  -13.93%  '

$ws.Range('D30').Value = '11.93'
$ws.Range('E30').Value = '  -4.28%  '

$ws.Range('D31').Value = '2.55'
$ws.Range('E31').Value = '  -8.79%  '

$ws.Range('E32').Value = '  -6.66%  '

$ws.Range('D33').Value = '6.80'
$ws.Range('E33').Value = '  -8.91%  '

$ws.Range('E34').Value = '  -8.77%  '

$ws.Range('E35').Value = '  +0.02%  '

$ws.Range('D36').Value = '36.75'
$ws.Range('E36').Value = '  -11.18%  '

$ws.Range('D37').Value = '53.59'
$ws.Range('E37').Value = '  -5.58%  '

$ws.Range('D38').Value = '0.0437'

$ws.Range('D39').Value = '0.996'
$ws.Range('E39').Value = '  -0.11%  '

$ws.Range('D40').Value = '2.78'
$ws.Range('E40').Value = '  +17.98%  '

$ws.Range('D41').Value = '0.0₃0633'
$ws.Range('E41').Value = '  -14.18%  '

$ws.Range('E42').Value = '  -10.77%  '

$ws.Range('D43').Value = '3.10'
$ws.Range('E43').Value = '  +13.17%  '

$ws.Range('D44').Value = '141.42'
$ws.Range('E44').Value = '  -4.85%  '

$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').Value = '25.28'
$ws.Range('E45').Value = '  +16.45%  '

$ws.Range('B46').Value = 'Stacks'
$ws.Range('C46').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D46').Value = '2.74'
$ws.Range('E46').Value = '  -10.31%  '

$ws.Range('D47').Value = '3.08'
$ws.Range('E47').Value = '  -6.37%  '

$ws.Range('D48').Value = '1.95'
$ws.Range('E48').Value = '  -2.28%  '

$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D49').Value = '4.03'
$ws.Range('E49').Value = '  -7.83%  '

$ws.Range('B50').Value = 'WEMIXToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D50').Value = '2.44'
$ws.Range('E50').Value = '  -10.97%  '

$ws.Range('D51').Value = '0.277'
$ws.Range('E51').Value = '  -11.06%  '
